$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate row 30 (same look/feel) into the new row 31, then edit its text.
$ws.Cells.Item(30, 1).Copy($ws.Cells.Item(31, 1))
$ws.Cells.Item(30, 2).Copy($ws.Cells.Item(31, 2))

$ws.Cells.Item(31, 1).Value = "hybrid_cbf_cf_w0.13cf_w0.87cbf_popularity1000_biasGiustiSkr5"

$fmt = $ws.Cells.Item(31, 2).NumberFormat
$ws.Cells.Item(31, 2).NumberFormat = "@"
$ws.Cells.Item(31, 2).Value = "0.00555"
$ws.Cells.Item(31, 2).NumberFormat = $fmt

# Mirror the author's navigation/selection after adding the row
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("A34").Select()
